# [Kadastro App] Kayıt silindi: 1
# Deletes the "Kayıt No" = 1 (Anamur / 18-UYG.) record from the "Kayitlar"
# sheet and the corresponding mirrored row on the "Anamur" filtered sheet,
# shifting the remaining rows up.

$wb = $excel.ActiveWorkbook

# Remove the record from the main "Kayitlar" list (row 3: Kayıt No = 1)
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(3).Delete()

# Remove the mirrored record from the "Anamur" filtered sheet (row 2)
$wsAnamur = $wb.Worksheets.Item("Anamur")
$wsAnamur.Rows.Item(2).Delete()
